# Update the "R30" rule row on the Rules sheet: the "From" value (C10)
# changes from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
